$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.60954197817034
$ws.Range("C2").Value = 6.393399622049738
$ws.Range("D2").Value = 6.484671421832314
$ws.Range("E2").Value = 16.30822471466247
$ws.Range("F2").Value = 35.06749170459229
$ws.Range("K2").Value = 12.91713939145572
$ws.Range("N2").Value = 20.891858391433
$ws.Range("B3").Value = 13.2992885757827
$ws.Range("C3").Value = 6.111244584545777
$ws.Range("D3").Value = 6.49869982611556
$ws.Range("E3").Value = 15.39321256051904
$ws.Range("F3").Value = 34.78736585981594
$ws.Range("K3").Value = 12.69698921044756
$ws.Range("N3").Value = 20.9195938942269
$ws.Range("B4").Value = 13.10965990068869
$ws.Range("C4").Value = 5.933454043404323
$ws.Range("D4").Value = 6.507556599173099
$ws.Range("E4").Value = 14.80825640629559
$ws.Range("F4").Value = 34.6231509947717
$ws.Range("K4").Value = 12.56406682283926
$ws.Range("N4").Value = 20.9384568268156
$ws.Range("B5").Value = 13.03272826527362
$ws.Range("C5").Value = 5.859983051374364
$ws.Range("D5").Value = 6.511227304035269
$ws.Range("E5").Value = 14.56434032838415
$ws.Range("F5").Value = 34.55824208636092
$ws.Range("K5").Value = 12.51054216047675
$ws.Range("N5").Value = 20.94660367206179
$ws.Range("B6").Value = 13.01997828563072
$ws.Range("C6").Value = 5.84772555155997
$ws.Range("D6").Value = 6.51184054558383
$ws.Range("E6").Value = 14.52351264597438
$ws.Range("F6").Value = 34.54758687092048
$ws.Range("K6").Value = 12.5016954376963
$ws.Range("N6").Value = 20.94798420787842
$ws.Range("B7").Value = 13.10862081530594
$ws.Range("C7").Value = 5.93246714334597
$ws.Range("D7").Value = 6.507605854198372
$ws.Range("E7").Value = 14.80498890344664
$ws.Range("F7").Value = 34.62226740577962
$ws.Range("K7").Value = 12.56334226974296
$ws.Range("N7").Value = 20.93856483664792
$ws.Range("B8").Value = 13.50246073389704
$ws.Range("C8").Value = 6.297125440170423
$ws.Range("D8").Value = 6.48945812067759
$ws.Range("E8").Value = 15.99768001179848
$ws.Range("F8").Value = 34.96932169457211
$ws.Range("K8").Value = 12.84081046580321
$ws.Range("N8").Value = 20.9010405474021
$ws.Range("B9").Value = 14.27632484620003
$ws.Range("C9").Value = 6.971478470939132
$ws.Range("D9").Value = 6.455785322062487
$ws.Range("E9").Value = 18.16694151412855
$ws.Range("F9").Value = 35.70915304109329
$ws.Range("K9").Value = 13.39941147107716
$ws.Range("N9").Value = 20.84204368229366
$ws.Range("B10").Value = 14.83898125040862
$ws.Range("C10").Value = 7.436798923953047
$ws.Range("D10").Value = 6.432190744138206
$ws.Range("E10").Value = 19.79230145112535
$ws.Range("F10").Value = 36.28525387245367
$ws.Range("K10").Value = 13.81416497874535
$ws.Range("N10").Value = 20.80765111702454
$ws.Range("B11").Value = 15.09232927512505
$ws.Range("C11").Value = 7.641072030997472
$ws.Range("D11").Value = 6.421700686622358
$ws.Range("E11").Value = 20.49057219916623
$ws.Range("F11").Value = 36.55358423532507
$ws.Range("K11").Value = 14.00286883827783
$ws.Range("N11").Value = 20.79396176333264
$ws.Range("B12").Value = 15.18778806604158
$ws.Range("C12").Value = 7.717301828672936
$ws.Range("D12").Value = 6.417762997475804
$ws.Range("E12").Value = 20.74911539155542
$ws.Range("F12").Value = 36.65602188190193
$ws.Range("K12").Value = 14.07425658488973
$ws.Range("N12").Value = 20.7890602798137
$ws.Range("B13").Value = 15.16725216986949
$ws.Range("C13").Value = 7.700935216203499
$ws.Range("D13").Value = 6.4186095123013
$ws.Range("E13").Value = 20.69369430052276
$ws.Range("F13").Value = 36.63392451363456
$ws.Range("K13").Value = 14.05888624086138
$ws.Range("N13").Value = 20.79010332551457
$ws.Range("B14").Value = 15.10019288827499
$ws.Range("C14").Value = 7.647366336390189
$ws.Range("D14").Value = 6.421376037533744
$ws.Range("E14").Value = 20.5119602811163
$ws.Range("F14").Value = 36.56199569761717
$ws.Range("K14").Value = 14.00874377750369
$ws.Range("N14").Value = 20.79355284973228
$ws.Range("B15").Value = 15.05905185548993
$ws.Range("C15").Value = 7.614405868361517
$ws.Range("D15").Value = 6.423075119774873
$ws.Range("E15").Value = 20.3998786606917
$ws.Range("F15").Value = 36.51804268920968
$ws.Range("K15").Value = 13.97801874830459
$ws.Range("N15").Value = 20.79570258965539
$ws.Range("B16").Value = 14.82236287387837
$ws.Range("C16").Value = 7.423294929176717
$ws.Range("D16").Value = 6.432881162112674
$ws.Range("E16").Value = 19.74584352617361
$ws.Range("F16").Value = 36.26783757199029
$ws.Range("K16").Value = 13.80182689957533
$ws.Range("N16").Value = 20.80858520284796
$ws.Range("B17").Value = 14.6764180950288
$ws.Range("C17").Value = 7.304113088791846
$ws.Range("D17").Value = 6.438958929172565
$ws.Range("E17").Value = 19.33410914472457
$ws.Range("F17").Value = 36.1158963162996
$ws.Range("K17").Value = 13.69369265663084
$ws.Range("N17").Value = 20.81699003489422
$ws.Range("B18").Value = 14.59223488412194
$ws.Range("C18").Value = 7.234868908005302
$ws.Range("D18").Value = 6.442477612752498
$ws.Range("E18").Value = 19.09341850779614
$ws.Range("F18").Value = 36.02909737193072
$ws.Range("K18").Value = 13.63150386155303
$ws.Range("N18").Value = 20.8220083135886
$ws.Range("B19").Value = 14.56369401187717
$ws.Range("C19").Value = 7.211306797032835
$ws.Range("D19").Value = 6.443672923258392
$ws.Range("E19").Value = 19.01125871118494
$ws.Range("F19").Value = 35.9998129081547
$ws.Range("K19").Value = 13.61045135656503
$ws.Range("N19").Value = 20.82373899202287
$ws.Range("B20").Value = 14.69197974060202
$ws.Range("C20").Value = 7.316872478736598
$ws.Range("D20").Value = 6.438309571436342
$ws.Range("E20").Value = 19.37833944310913
$ws.Range("F20").Value = 36.13200980873474
$ws.Range("K20").Value = 13.70520354053444
$ws.Range("N20").Value = 20.81607627072454
$ws.Range("B21").Value = 15.11990360562713
$ws.Range("C21").Value = 7.663131747590363
$ws.Range("D21").Value = 6.420562503579961
$ws.Range("E21").Value = 20.56549918170506
$ws.Range("F21").Value = 36.58310108768826
$ws.Range("K21").Value = 14.02347432312989
$ws.Range("N21").Value = 20.79253196837691
$ws.Range("B22").Value = 15.39674379772453
$ws.Range("C22").Value = 7.882859397866603
$ws.Range("D22").Value = 6.409165698470227
$ws.Range("E22").Value = 21.30713855461095
$ws.Range("F22").Value = 36.88269940358328
$ws.Range("K22").Value = 14.2310403880351
$ws.Range("N22").Value = 20.77879083166588
$ws.Range("B23").Value = 15.24928061014623
$ws.Range("C23").Value = 7.766205323707713
$ws.Range("D23").Value = 6.415230012740222
$ws.Range("E23").Value = 20.91443283722317
$ws.Range("F23").Value = 36.72238496427213
$ws.Range("K23").Value = 14.12032283203506
$ws.Range("N23").Value = 20.78597372117634
$ws.Range("B24").Value = 14.68494517735676
$ws.Range("C24").Value = 7.311106212038965
$ws.Range("D24").Value = 6.438603069710258
$ws.Range("E24").Value = 19.35835532231429
$ws.Range("F24").Value = 36.12472316713601
$ws.Range("K24").Value = 13.69999952952016
$ws.Range("N24").Value = 20.81648880343244
$ws.Range("B25").Value = 14.06753569357259
$ws.Range("C25").Value = 6.794005497978477
$ws.Range("D25").Value = 6.46469204244773
$ws.Range("E25").Value = 17.58471719701079
$ws.Range("F25").Value = 35.50302250119584
$ws.Range("K25").Value = 13.24721134675372
$ws.Range("N25").Value = 20.85643602709017
